# Use Case cell originally reads "UI depicting Enemy States".
# The edit capitalizes the "d" in "depicting" -> "Depicting", which (per the
# canonical OOXML diff) is realized as the single run being split into three
# runs: "UI ", "D", "epicting Enemy States" (as if the lower-case "d" had
# been selected and retyped as an upper-case "D").

$d = $word.ActiveDocument

# Locate the exact phrase so we don't touch any other cell in the document.
$target = $d.Content
$found = $target.Find.Execute("UI depicting Enemy States", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find 'UI depicting Enemy States' in the document"
}

$phraseStart = $target.Start

# The "d" that needs to become "D" is the 4th character of the phrase
# ("UI " is 3 characters long).
$dStart = $phraseStart + 3
$dChar = $d.Range($dStart, $dStart + 1)

# Briefly toggle Bold on just that character: changing its direct formatting
# forces Word to split the run at this position (matching the run boundaries
# introduced by the real edit), then type the capital "D" over it.
$dChar.Font.Bold = $true
$dChar.Text = "D"

# Restore the original (non-bold) formatting so the visible/effective
# formatting of the paragraph is unchanged - only the run is left split.
$dCharAgain = $d.Range($dStart, $dStart + 1)
$dCharAgain.Font.Bold = $false
